$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.108.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.29%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.506.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.14%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "419.14"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.00"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.74%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.654"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.25%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.779"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.80%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +13.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.68"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.98%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000262"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +19.64%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.10"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +9.70%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.060.52"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.79%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.53"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.03%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.493.76"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.35%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.80"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.73%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.11"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.06%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "65.010.85"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.09%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "455.73"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.93%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.24"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.75%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.56%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.37"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.97"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.22%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "34.11"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.76%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.59"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.07%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.73"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.52"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.49%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.118"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.56%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "40.13"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.16%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.05"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.91%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0509"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.03%  "

# Row 37
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0740"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +34.63%  "

# Row 38
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.150"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.48%  "

# Row 39
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.18%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.07"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.38%  "

# Row 41
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.57"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.37%  "

# Row 42
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.29"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.58%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.31"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.53%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.315"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.55%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.26%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.33"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.73%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.92"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.87%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.145"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.25%  "

# Row 50
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.58"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +12.15%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.85"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.57%  "
